$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: fill column A (name) for new rows in the order strings were originally entered
$ws.Range("A5").Value = "Chicken"
$ws.Range("A6").Value = "Pork"
$ws.Range("A9").Value = "Cheese"
$ws.Range("A10").Value = "Lard"
$ws.Range("A15").Value = "Milkshakes"
$ws.Range("A11").Value = "Greens"
$ws.Range("A12").Value = "Vegetables"
$ws.Range("A13").Value = "Fruits"
$ws.Range("A14").Value = "Drinks"
$ws.Range("A16").Value = "Soda"
$ws.Range("A17").Value = "Tea"
$ws.Range("A18").Value = "Baking"
$ws.Range("A19").Value = "Bread"
$ws.Range("A20").Value = "Cookies"
$ws.Range("A21").Value = "Pies"

# Step 2: fill column B (parent_string) for new/shifted rows (reuses existing strings)
$ws.Range("B5").Value = "Meat"
$ws.Range("B6").Value = "Meat"
$ws.Range("B9").Value = "Dairy"
$ws.Range("B10").Value = "Dairy"
$ws.Range("B11").Value = "Food"
$ws.Range("B12").Value = "Greens"
$ws.Range("B13").Value = "Greens"
$ws.Range("B14").Value = "Food"
$ws.Range("B15").Value = "Drinks"
$ws.Range("B16").Value = "Drinks"
$ws.Range("B17").Value = "Drinks"
$ws.Range("B18").Value = "Food"
$ws.Range("B19").Value = "Baking"
$ws.Range("B20").Value = "Baking"
$ws.Range("B21").Value = "Baking"

# Step 3: fill column C (picture) for new rows in the same order as column A
$ws.Range("C5").Value = "Chicken.jfif"
$ws.Range("C6").Value = "Pork.jpg"
$ws.Range("C9").Value = "Cheese.jpg"
$ws.Range("C10").Value = "Lard.jpg"
$ws.Range("C15").Value = "Milkshakes.jpg"
$ws.Range("C11").Value = "Greens.jpg"
$ws.Range("C12").Value = "Vegetables.jpg"
$ws.Range("C13").Value = "Fruits.jpg"
$ws.Range("C14").Value = "Drinks.jpg"
$ws.Range("C16").Value = "Soda.jpg"
$ws.Range("C17").Value = "Tea.jpg"
$ws.Range("C18").Value = "Baking.jpg"
$ws.Range("C19").Value = "Bread.jpg"
$ws.Range("C20").Value = "Cookies.jpg"
$ws.Range("C21").Value = "Pies.jpg"

# Step 4: restore shifted Dairy/Milk rows (7,8) - same content as originally rows 5,6
$ws.Range("A7").Value = "Dairy"
$ws.Range("B7").Value = "Food"
$ws.Range("C7").Value = "Dairy.jfif"
$ws.Range("A8").Value = "Milk"
$ws.Range("B8").Value = "Dairy"
$ws.Range("C8").Value = "Milk.jfif"
$ws.Range("B15").Select()
